$wb = $excel.ActiveWorkbook

# The "Granty_przyznane" sheet (sheet7) contains three yearly subtotal rows
# (one per year: 2021, 2020, 2019) whose "Jednostka" column reads "Razem"
# (i.e. "Total"). These duplicate the totals already shown on the "Razem"
# sheet, so they are being removed from the per-faculty grants listing.
$wsGranted = $wb.Worksheets.Item("Granty_przyznane")

# Walk the used rows bottom-up (so deleting a row doesn't shift the index of
# rows we still need to inspect) and remove every row whose first column is
# the literal text "Razem".
$lastRow = $wsGranted.UsedRange.Rows.Count
for ($r = $lastRow; $r -ge 1; $r--) {
    $label = $wsGranted.Cells.Item($r, 1).Value2
    if ($label -eq "Razem") {
        $wsGranted.Rows($r).EntireRow.Delete()
    }
}

# The edit session ends with the "Granty_przyznane" tab active (previously
# "Granty_złożone" was the selected/active tab) ...
$wsGranted.Activate() | Out-Null

# ... with the last data row selected as a full-row selection, matching the
# state captured after scrolling through / deleting the trailing rows.
$lastDataRow = $wsGranted.UsedRange.Rows.Count + 1
$wsGranted.Range("A" + $lastDataRow + ":XFD" + $lastDataRow).Select() | Out-Null

# Scroll the view so row 16 is at the top (best-effort: not all hosts persist
# window scroll position, but this matches the recorded view state).
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
